$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column C for rows 2 through 108
# from serial date 45233 (2023-11-03) to 45243 (2023-11-13)
$ws.Range("C2:C108").Value = 45243
